$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.951.39'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").Value = '2.933.53'
$ws.Range("E3").Value = '  -1.98%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''373.01'
$ws.Range("E5").Value = '  -2.10%  '

$ws.Range("D6").Value = '''101.98'
$ws.Range("E6").Value = '  -3.93%  '

$ws.Range("D7").Value = '''0.535'
$ws.Range("E7").Value = '  -2.01%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = '''0.581'
$ws.Range("E9").Value = '  -3.28%  '

$ws.Range("D10").Value = '''36.48'
$ws.Range("E10").Value = '  -2.75%  '

$ws.Range("E11").Value = '  -0.81%  '

$ws.Range("D12").Value = '''0.0833'
$ws.Range("E12").Value = '  -1.63%  '

$ws.Range("D13").Value = '3.398.99'
$ws.Range("E13").Value = '  -2.10%  '

$ws.Range("D14").Value = '''17.85'
$ws.Range("E14").Value = '  -4.33%  '

$ws.Range("D15").Value = '''7.32'
$ws.Range("E15").Value = '  -2.89%  '

$ws.Range("D16").Value = '2.939.47'
$ws.Range("E16").Value = '  -1.92%  '

$ws.Range("D17").Value = '''0.973'
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").Value = '50.983.66'
$ws.Range("E18").Value = '  -1.37%  '

$ws.Range("D19").Value = '''3.14'
$ws.Range("E19").Value = '  -7.32%  '

$ws.Range("E20").Value = '  -4.52%  '

$ws.Range("D21").Value = '''12.52'
$ws.Range("E21").Value = '  -4.28%  '

$ws.Range("D22").Value = '0.0₃0951'
$ws.Range("E22").Value = '  -1.40%  '

$ws.Range("D23").Value = '''263.09'
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").Value = '''68.12'
$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("E25").Value = '  +2.00%  '

$ws.Range("D26").Value = '''7.87'
$ws.Range("E26").Value = '  +8.71%  '

$ws.Range("D27").Value = '''8.03'
$ws.Range("E27").Value = '  +6.98%  '

$ws.Range("E28").Value = '  -1.70%  '

$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.113'
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").Value = '''25.61'
$ws.Range("E31").Value = '  -1.90%  '

$ws.Range("D32").Value = '''9.82'
$ws.Range("E32").Value = '  -0.77%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '''50.99'
$ws.Range("E33").Value = '  -0.71%  '

$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").Value = '''0.0455'
$ws.Range("E34").Value = '  -1.44%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").Value = '''33.70'
$ws.Range("E35").Value = '  -3.97%  '

$ws.Range("D36").Value = '''2.02'
$ws.Range("E36").Value = '  -3.48%  '

$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").Value = '''2.97'
$ws.Range("E38").Value = '  -4.72%  '

$ws.Range("D39").Value = '''2.53'
$ws.Range("E39").Value = '  -3.25%  '

$ws.Range("D40").Value = '''0.114'
$ws.Range("E40").Value = '  -1.71%  '

$ws.Range("D41").Value = '''16.33'
$ws.Range("E41").Value = '  -7.13%  '

$ws.Range("D42").Value = '''1.77'
$ws.Range("E42").Value = '  -4.82%  '

$ws.Range("D43").Value = '''120.80'
$ws.Range("E43").Value = '  -2.95%  '

$ws.Range("D44").Value = '''20.92'
$ws.Range("E44").Value = '  -6.89%  '

$ws.Range("E45").Value = '  -1.66%  '

$ws.Range("D46").Value = '''0.274'
$ws.Range("E46").Value = '  -2.19%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '''2.30'
$ws.Range("E47").Value = '  -4.58%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '''3.20'
$ws.Range("E48").Value = '  -1.98%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.987.46'
$ws.Range("E49").Value = '  -3.29%  '

$ws.Range("D50").Value = '''0.0344'
$ws.Range("E50").Value = '  -5.63%  '

$ws.Range("D51").Value = '''5.02'
$ws.Range("E51").Value = '  -3.45%  '
